$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B width: widen to fit the new, longer item descriptions ---
$ws.Columns("B").ColumnWidth = 33.7

# --- Row 12: Điện trở 2K 0805 + Tụ 10uF 0805 / thegioiic.com ---
$ws.Range("A9:D9").Copy()
$ws.Range("A12:D12").PasteSpecial(-4122)
$ws.Range("B12").Value = "Điện trở 2K 0805 + Tụ 10uF 0805"
$ws.Range("D12").Value = "http://thegioiic.com/"
$ws.Range("A12").Value = (Get-Date -Year 2018 -Month 5 -Day 18 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C12").Value = 14800

# --- Row 11: Diode 1N4007 + Nút + Lưỡi cưa mạch + Bút lông / Bảo tín ---
$ws.Range("A9:D9").Copy()
$ws.Range("A11:D11").PasteSpecial(-4122)
$ws.Range("D11").Value = "Bảo tín"
$ws.Range("B11").Value = "Diode 1N4007 + Nút + Lưỡi cưa mạch + Bút lông"
$ws.Range("A11").Value = (Get-Date -Year 2018 -Month 5 -Day 17 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C11").Value = 65000

# --- Row 10: GPS Ublox Neo7, linked to hshop.vn ---
$ws.Range("A9:D9").Copy()
$ws.Range("A10:D10").PasteSpecial(-4122)
$ws.Range("B10").Value = "GPS Ublox Neo7"
$ws.Range("D10").Value = "http://hshop.vn/products/mach-gps-ublox-neo-7"
$ws.Range("C10").Value = 380000
$ws.Range("D9").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("D10"), "http://hshop.vn/products/mach-gps-ublox-neo-7")

$ws.Range("D21").Select()
